$d = $word.ActiveDocument

# Locate the unique target sentence (the 3rd "Sub goal..." paragraph, about finger counting).
$target = "Sub goal are counting to 10, counting to 100, and counting to 1000, then figuring out what finger she stops on."
$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target sentence"
}
$matchStart = $rng.Start
$matchEnd = $rng.End

# --- Step 1: insert "s" after "goal" (before the space), turning "Sub goal" into "Sub goals".
$posAfterGoal = $matchStart + 8
$insPt = $d.Range($posAfterGoal, $posAfterGoal)
$insPt.InsertAfter("s")

# --- Step 2: append the new trailing sentence (after "stops on.") with matching
# Times New Roman formatting.
$newEnd = $matchEnd + 1
$endPt = $d.Range($newEnd, $newEnd)
$newSentence = " Possible solutions for sub goals are counting to ten then counting to on hundred and find out what fingers you stop on."
$endPt.InsertAfter($newSentence)

# --- Step 3: now that all text content is final, restore/introduce the run
# boundaries the live edit session would have produced. A transient
# bookmark forces a split at a position; deleting the bookmark right after
# leaves the split in place without leaving stray markers behind.
function Split-At($pos) {
    $p = $d.Range($pos, $pos)
    $d.Bookmarks.Add("ZZZSPLITMARK", $p)
    $d.Bookmarks("ZZZSPLITMARK").Delete()
}

Split-At ($matchStart - 1)          # boundary before " Sub goal..." (was its own run originally)
Split-At $posAfterGoal               # boundary between "Sub goal" and "s"
$posAfterS = $posAfterGoal + 1
Split-At $posAfterS                  # boundary right after "s" (where _GoBack will sit)

$thenStart = $matchStart + 1 + 68    # start of "then" (shifted by 1 for the inserted "s")
$thenEnd = $thenStart + 4            # end of "then"
Split-At $thenStart
Split-At $thenEnd

Split-At $newEnd                     # boundary before the newly appended sentence

# --- Step 4: move the _GoBack bookmark to sit right after the new "s"
# (before " are counting..."), matching where the live edit cursor last was.
$goBackPt = $d.Range($posAfterS, $posAfterS)
$d.Bookmarks.Add("_GoBack", $goBackPt)

# --- Step 5: make sure the newly appended sentence carries the same
# Times New Roman font as the rest of the paragraph.
$addedRange = $d.Range($newEnd, $newEnd + $newSentence.Length)
$addedRange.Font.Name = "Times New Roman"

Write-Output "done"
